# Update "想去人数" (want-to-go count) figures in column F for the
# 展览 (Exhibitions) sheet and the aggregated 全部类型 (All types) sheet.
# These are plain numeric value updates, no formulas or formatting involved.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 149
$wsExpo.Range("F3").Value = 40
$wsExpo.Range("F4").Value = 235
$wsExpo.Range("F5").Value = 3845
$wsExpo.Range("F6").Value = 25

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 149
$wsAll.Range("F3").Value = 40
$wsAll.Range("F4").Value = 235
$wsAll.Range("F5").Value = 3845
$wsAll.Range("F8").Value = 25
